$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1) Rows 31-50: flip sign on A/D/E/F and bump C from 2 -> 3
#    (B and G stay as-is: 4000 and 0)
# ---------------------------------------------------------------------------
for ($r = 31; $r -le 50; $r++) {
    $ws.Cells.Item($r, 1).Value2 = -1000   # A
    $ws.Cells.Item($r, 3).Value2 = 3       # C
    $ws.Cells.Item($r, 4).Value2 = -50     # D
    $ws.Cells.Item($r, 5).Value2 = -500    # E
    $ws.Cells.Item($r, 6).Value2 = -150    # F
}

# ---------------------------------------------------------------------------
# 2) Append new rows 61-100 repeating the original 1000/4000/2/50/500/150/0
#    pattern that used to occupy rows 31-60.
# ---------------------------------------------------------------------------
for ($r = 61; $r -le 100; $r++) {
    $ws.Cells.Item($r, 1).Value2 = 1000    # A
    $ws.Cells.Item($r, 2).Value2 = 4000    # B
    $ws.Cells.Item($r, 3).Value2 = 2       # C
    $ws.Cells.Item($r, 4).Value2 = 50      # D
    $ws.Cells.Item($r, 5).Value2 = 500     # E
    $ws.Cells.Item($r, 6).Value2 = 150     # F
    $ws.Cells.Item($r, 7).Value2 = 0       # G
}

# ---------------------------------------------------------------------------
# 3) Update the view state: scroll so row 65 is at the top and select
#    A30:G50 with A30 as the active cell.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 65
[void]$ws.Range("A30:G50").Select()
